$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17: new entries "waypoints" / "read waypoints"
$ws.Range("C17").Value = "waypoints"
$ws.Range("F17").Value = "read waypoints"

# Row 18: new entries "camera values" / "read camera values"
$ws.Range("C18").Value = "camera values"
$ws.Range("F18").Value = "read camera values"

# Row 9: "transmitter calibration" -> "transmitter slope cal"
#        "read transmitter calibration" -> "read transmitter slope values"
$ws.Range("C9").Value = "transmitter slope cal"

# Row 10: "camera values" -> "transmitter offset cal"
#         "read camera values" -> "read transmitter offset values"
$ws.Range("C10").Value = "transmitter offset cal"

$ws.Range("F9").Value = "read transmitter slope values"
$ws.Range("F10").Value = "read transmitter offset values"

# Update selected cell to F10
$ws.Range("F10").Select()
